$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:H1) ---------------------------------------------------
# Extend the existing bold/centered/bordered header style (A1:G1) onto the
# new H1 header cell by copying formats from G1 before writing values.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("A1").Value = "Employee"
$ws.Range("B1").Value = "Department"
$ws.Range("C1").Value = "Salary"
$ws.Range("D1").Value = "Bonus"
$ws.Range("E1").Value = "Performance"
$ws.Range("F1").Value = "Quota Met"
$ws.Range("G1").Value = "Start Date"
$ws.Range("H1").Value = "Hours Worked"

# --- Data rows (A2:H6) -----------------------------------------------------
$ws.Range("A2").Value = "Alice Johnson"
$ws.Range("B2").Value = "Engineering"
$ws.Range("C2").Value = 85000
$ws.Range("D2").Value = 8500
$ws.Range("E2").Value = 0.92
$ws.Range("F2").Value = 0.75
$ws.Range("G2").Value = 44270
$ws.Range("H2").Value = 42.5

$ws.Range("A3").Value = "Bob Smith"
$ws.Range("B3").Value = "Sales"
$ws.Range("C3").Value = 72000
$ws.Range("D3").Value = 10800
$ws.Range("E3").Value = 0.85
$ws.Range("F3").Value = 1.15
$ws.Range("G3").Value = 43668
$ws.Range("H3").Value = 38.75

$ws.Range("A4").Value = "Carol White"
$ws.Range("B4").Value = "Marketing"
$ws.Range("C4").Value = 68000
$ws.Range("D4").Value = 5100
$ws.Range("E4").Value = 0.78
$ws.Range("F4").Value = 0.92
$ws.Range("G4").Value = 44571
$ws.Range("H4").Value = 40

$ws.Range("A5").Value = "David Brown"
$ws.Range("B5").Value = "Engineering"
$ws.Range("C5").Value = 92000
$ws.Range("D5").Value = 11040
$ws.Range("E5").Value = 0.95
$ws.Range("F5").Value = 0.88
$ws.Range("G5").Value = 43409
$ws.Range("H5").Value = 45.25

$ws.Range("A6").Value = "Eva Martinez"
$ws.Range("B6").Value = "Sales"
$ws.Range("C6").Value = 78000
$ws.Range("D6").Value = 9360
$ws.Range("E6").Value = 0.88
$ws.Range("F6").Value = 1.05
$ws.Range("G6").Value = 44000
$ws.Range("H6").Value = 39.5

# --- Number formats (applied to the data rows only) ------------------------
# Stamp the full date+time format onto G2 first (mirrors the author's
# original formatting pass) before narrowing the whole Start Date column
# down to a date-only display.
$ws.Range("G2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("C2:D6").NumberFormat = "$#,##0"
$ws.Range("E2:F6").NumberFormat = "0%"
$ws.Range("G2:G6").NumberFormat = "YYYY-MM-DD"
$ws.Range("H2:H6").NumberFormat = "0.00"
